$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.886.59"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "3.845.11"
$ws.Range("E3").Value = "  +5.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "422.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("D7").Value = "3.839.23"
$ws.Range("E7").Value = "  +5.65%  "
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.720"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  +1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.75%  "
$ws.Range("D15").Value = "4.451.90"
$ws.Range("E15").Value = "  +5.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +14.75%  "
$ws.Range("D17").Value = "3.862.28"
$ws.Range("E17").Value = "  +6.19%  "
$ws.Range("E18").Value = "  -0.70%  "
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "67.105.35"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "409.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.63%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +35.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "746.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.121"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.78%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -5.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("E39").Value = "  +22.96%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").Value = "0.0₃0723"
$ws.Range("E41").Value = "  +10.45%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("E44").Value = "  +1.45%  "
$ws.Range("E45").Value = "  -4.63%  "
$ws.Range("E46").Value = "  +8.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("E51").Value = "  +0.38%  "
